# Applies the "NIT-9010645851" update:
#  - Adds a new worker (MAYELIS MARIA DEL TORO MANGONES) as a new first data row.
#  - Re-orders the "Periodo Mora" rows for CARLOS FERNANDO PUENTES LIZARAZO.
#  - Refreshes the totals (Valor Mora, Cant. Trabajadores, Cant. Periodos).
#  - Widens a few columns to fit the new/longer content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new data row right above the existing first worker row (16) ---
# Shift cells down, pulling formatting from the row that is about to become
# row 17 (the previous row 16) so the new row re-uses the existing "inner
# row" style (borders/fill/number formats) instead of Excel minting new ones.
$ws.Rows.Item(16).Insert(-4121, 0)
$ws.Range("B17:J17").Copy() | Out-Null
$ws.Range("B16:J16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- 2. Populate the new worker row ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1067164354"
$ws.Range("D16").Value = "MAYELIS MARIA DEL TORO MANGONES"
$ws.Range("E16").Value = "2503"
$ws.Range("F16").Value = 1898
$ws.Range("G16").Value = 1423500

# --- 3. Row 17 (GUIDO BALDIRIS ESCORCIA) keeps its values; nothing to change ---

# --- 4. Re-order the "Periodo Mora" entries for CARLOS FERNANDO PUENTES LIZARAZO ---
# Row 18 -> period 1910 / valor mora 44000
$ws.Range("E18").Value = "1910"
$ws.Range("F18").Value = 44000

# Row 19 -> period 1909 / valor mora 44000 (unchanged, set explicitly for safety)
$ws.Range("E19").Value = "1909"
$ws.Range("F19").Value = 44000

# Row 20 -> period 1908 / valor mora 11733
$ws.Range("E20").Value = "1908"
$ws.Range("F20").Value = 11733

# --- 5. Refresh header totals ---
$ws.Range("E11").Value = 107881
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 5

# --- 6. Widen columns to better fit the new/longer content ---
$ws.Columns.Item(2).ColumnWidth = 17.7
$ws.Columns.Item(3).ColumnWidth = 15.9
$ws.Columns.Item(5).ColumnWidth = 12.7
$ws.Columns.Item(6).ColumnWidth = 9.3
$ws.Columns.Item(7).ColumnWidth = 13.5
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 17.3
$ws.Columns.Item(10).ColumnWidth = 14.2

Write-Output "Applied NIT-9010645851 update"
